# Applies the 6-Jan-2023 GitHub Actions "symbol list" refresh to the crypto
# price sheet: price (D), 1h volume % (E) updates on most rows, plus a
# coin-name/link (B/C) re-rank shuffle across rows 9-15.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without letting Excel
# auto-convert numeric-/percent-looking strings into numbers (which would
# reformat "4.650" -> "4.65" or "10,438.14%" -> a recomputed percent).
function Set-TextValue([string]$addr, [string]$text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '255.87'
Set-TextValue 'E2' '-0.33%'

# Row 3
Set-TextValue 'D3' '26.43'
Set-TextValue 'E3' '-2.44%'

# Row 4
Set-TextValue 'D4' '4.650'
Set-TextValue 'E4' '-0.58%'

# Row 5
Set-TextValue 'D5' '0.05924'
Set-TextValue 'E5' '0.14%'

# Row 6
Set-TextValue 'D6' '6.614'
Set-TextValue 'E6' '-0.41%'

# Row 7
Set-TextValue 'D7' '0.8517'
Set-TextValue 'E7' '-1.91%'

# Row 8
Set-TextValue 'D8' '0.9092'
Set-TextValue 'E8' '-4.24%'

# Row 9
Set-TextValue 'B9' 'WazirX'
Set-TextValue 'C9' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D9' '0.1377'
Set-TextValue 'E9' '-2.27%'

# Row 10
Set-TextValue 'B10' 'LiechtensteinCryptoassetsExchange'
Set-TextValue 'C10' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D10' '0.04133'
Set-TextValue 'E10' '9.95%'

# Row 11
Set-TextValue 'B11' 'MandalaExchangeToken'
Set-TextValue 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D11' '0.06991'
Set-TextValue 'E11' '-1.43%'

# Row 12
Set-TextValue 'B12' 'BitrueCoin'
Set-TextValue 'C12' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D12' '0.03041'
Set-TextValue 'E12' '-4.89%'

# Row 13
Set-TextValue 'B13' 'BitMartToken'
Set-TextValue 'C13' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D13' '0.09086'
Set-TextValue 'E13' '-2.00%'

# Row 14
Set-TextValue 'B14' 'BitForexToken'
Set-TextValue 'C14' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D14' '0.001532'
Set-TextValue 'E14' '-0.73%'

# Row 15
Set-TextValue 'B15' 'One'
Set-TextValue 'C15' 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue 'D15' '0.0006029'
Set-TextValue 'E15' '-94.22%'

# Row 16
Set-TextValue 'D16' '0.006052'
Set-TextValue 'E16' '-0.57%'

# Row 18
Set-TextValue 'D18' '3.150'
Set-TextValue 'E18' '-1.37%'

# Row 19
Set-TextValue 'E19' '-1.91%'

# Row 20
Set-TextValue 'E20' '-3.58%'

# Row 21
Set-TextValue 'D21' '0.1285'
Set-TextValue 'E21' '0.09%'

# Row 22
Set-TextValue 'D22' '3.856'
Set-TextValue 'E22' '-0.18%'

# Row 23
Set-TextValue 'D23' '0.04218'
Set-TextValue 'E23' '-0.54%'

# Row 24
Set-TextValue 'D24' '0.001213'
Set-TextValue 'E24' '-0.82%'

# Row 25
Set-TextValue 'D25' '0.004685'
Set-TextValue 'E25' '9.08%'

# Row 26
Set-TextValue 'D26' '0.0001199'
Set-TextValue 'E26' '-0.11%'

# Row 40
Set-TextValue 'D40' '0.03783'
Set-TextValue 'E40' '-0.95%'

# Row 41
Set-TextValue 'D41' '0.006212'
Set-TextValue 'E41' '-0.59%'

# Row 42
Set-TextValue 'D42' '0.1095'
Set-TextValue 'E42' '-0.52%'

# Row 43
Set-TextValue 'D43' '0.002309'
Set-TextValue 'E43' '4.94%'

# Row 44
Set-TextValue 'D44' '0.01401'
Set-TextValue 'E44' '22.33%'

# Row 45
Set-TextValue 'E45' '-6.31%'

# Row 46
Set-TextValue 'D46' '0.00000000749'
Set-TextValue 'E46' '-0.12%'

# Row 47
Set-TextValue 'D47' '0.04999'
Set-TextValue 'E47' '-37.90%'

# Row 48
Set-TextValue 'E48' '10,438.14%'

# Row 49
Set-TextValue 'D49' '0.00002098'
Set-TextValue 'E49' '-0.12%'

# Row 50
Set-TextValue 'D50' '0.0001998'
Set-TextValue 'E50' '-0.12%'
